$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Refresh the cached "datetimeFigureOut" footer field (6/12/2014 -> 7/24/2014)
#    on the slide master and every slide layout.
# ---------------------------------------------------------------------------
function Set-DatePlaceholderText($shapes, $text) {
  for ($j = 1; $j -le $shapes.Count; $j++) {
    $sh = $shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
      $sh.TextFrame.TextRange.Text = $text
    }
  }
}

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes "7/24/2014"
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
  Set-DatePlaceholderText $master.CustomLayouts.Item($i).Shapes "7/24/2014"
}

# ---------------------------------------------------------------------------
# 2. Append two new "Title and Content" slides at the end of the deck.
# ---------------------------------------------------------------------------
$count = $p.Slides.Count

$s16 = $p.Slides.Add($count + 1, 2)
$s16.Shapes.Item(1).TextFrame.TextRange.Text = "User Stories"
$s16.Shapes.Item(2).TextFrame.TextRange.Text = "Wiki Page`r"

$s17 = $p.Slides.Add($count + 2, 2)
$s17.Shapes.Item(2).TextFrame.TextRange.Text = "Post page"
